# Introduce and populate a "Birthdate" date column (G) on the Person sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell + date-formatted body cells for the new column.
$ws.Range("G1:G3").NumberFormat = "dd/mm/yy"

$ws.Range("G1").Value = "Birthdate"
$ws.Range("G2").Value = 33193   # 16/11/90
$ws.Range("G3").Value = 33955   # 17/12/92

# Re-apply the (already Arial) font to the ZIP column so it carries an
# explicit font reference, matching the style used for the new column.
$ws.Range("D1:D3").Font.Name = "Arial"
